# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.546.91"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.43%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.246.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -4.53%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "494.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.02%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.12"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.24%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.24%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.78%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.291.11"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.08%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.81%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.44%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.83%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.651.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.39%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.02%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "54.380.86"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.67%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.91%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.253.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.49%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.28%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "304.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.30%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.49%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.46%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.01"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.65%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.374"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.30%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.42%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.362.34"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.65%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.49%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.01"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.93%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.57%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0686"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.60%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.23%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.993"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.21%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.50%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.59"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.50%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.67%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.865"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.44%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.36%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.07%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.28%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.64%  "

# Row 45
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.23%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "128.57"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.79%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.84"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.73%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0896"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.551"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "242.47"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.25%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.30%  "
